$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(1.505614041169197, 1.65323645889881, 16.98373111632243, 6.48142807727062, 26.62400969366105)
    3  = @(0.1554434735375247, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 1.705647867635037)
    4  = @(3.182878228561681, 0.3375848360084654, 0.1529057820181812, 0.4998867070740569, 4.173255553662385)
    5  = @(1.505614041169197, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 3.055818435266709)
    6  = @(0.02258322285507441, 0.05231270169004087, 0.1529057820181812, 0.4998867070740569, 0.7276884136373534)
    7  = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 6.48142807727062, 14.40014219143469)
    8  = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 8.418600821238126)
    9  = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 8.418600821238126)
    10 = @(1.505614041169197, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 6.741336633845642)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B - TB
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C - d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D - K
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E - IP
    $ws.Cells.Item($row, 7).Value = $vals[4]  # G - sum
}
